# Update investment cost results with newer values from server run.
$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 9739.537847600008
$ws.Range("E2").Value = 289823.7596598056
$ws.Range("I2").Value = 161752.8135478
$ws.Range("L2").Value = 485245.29503538
$ws.Range("M2").Value = 105905.87968015
$ws.Range("N2").Value = 70831.955579581
$ws.Range("O2").Value = 69610.4422391004

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 47386.06393082884
$ws.Range("E2").Value = 271236.7992183856
$ws.Range("I2").Value = 280426.171173861
$ws.Range("L2").Value = 184420.4799505123
$ws.Range("M2").Value = 113936.92264746
$ws.Range("N2").Value = 33931.8246116005
$ws.Range("O2").Value = 50485.47232467777

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 28619.61401238371
$ws.Range("B2").Value = 23143.29485244409
$ws.Range("E2").Value = 111916.8406725409
$ws.Range("I2").Value = 150385.2728707001
$ws.Range("M2").Value = 34803.41203795493
$ws.Range("N2").Value = 44938.11408779013
$ws.Range("O2").Value = 26938.31306104351

# --- Sheet "2040" ---
$ws = $wb.Worksheets.Item("2040")
$ws.Range("N2").Value = 1014.766490779938

# --- Sheet "2045" ---
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 34409.11717595647
$ws.Range("N2").Value = 5182.698656944208
$ws.Range("O2").Value = 22972.54525065906
